$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")
# ALC row 2
$ws_ALC.Range("H2").Value = 70
$ws_ALC.Range("I2").Value = 70
$ws_ALC.Range("J2").Value = 0
$ws_ALC.Range("K2").Value = 70
$ws_ALC.Range("L2").Value = 0
$ws_ALC.Range("M2").Value = 43
$ws_ALC.Range("N2").ClearContents() | Out-Null

# ALC row 19
$ws_ALC.Range("H19").Value = 1631.5883
$ws_ALC.Range("I19").Value = 1948.6364
$ws_ALC.Range("J19").Value = 1050.3334
$ws_ALC.Range("K19").Value = 1948.6364
$ws_ALC.Range("L19").Value = 1050.3334
$ws_ALC.Range("M19").Value = -1773.6364
$ws_ALC.Range("N19").Value = -1400.3334

# ALC row 74
$ws_ALC.Range("H74").Value = 5066.375
$ws_ALC.Range("J74").Value = 6588.154
$ws_ALC.Range("L74").Value = 6588.154
$ws_ALC.Range("N74").Value = -8460.154

# ALC row 77
$ws_ALC.Range("H77").Value = 5066.375
$ws_ALC.Range("J77").Value = 6588.154
$ws_ALC.Range("L77").Value = 32940.77
$ws_ALC.Range("N77").Value = -42300.77

# ALC row 107
$ws_ALC.Range("H107").Value = 1043.5555
$ws_ALC.Range("I107").Value = 986.5
$ws_ALC.Range("K107").Value = 986.5
$ws_ALC.Range("M107").Value = 933.5

# ALC row 125
$ws_ALC.Range("H125").Value = 1090.1052
$ws_ALC.Range("I125").Value = 1486
$ws_ALC.Range("J125").Value = 802.1818
$ws_ALC.Range("K125").Value = 13374
$ws_ALC.Range("L125").Value = 7219.6362
$ws_ALC.Range("M125").Value = -10914
$ws_ALC.Range("N125").Value = -12139.6362

# ALC row 132
$ws_ALC.Range("H132").Value = 49024.26
$ws_ALC.Range("I132").Value = 56058.535
$ws_ALC.Range("K132").Value = 168175.605
$ws_ALC.Range("M132").Value = -165645.605

# ALC row 138
$ws_ALC.Range("H138").Value = 5550.9707
$ws_ALC.Range("I138").Value = 7049.636
$ws_ALC.Range("K138").Value = 21148.908
$ws_ALC.Range("M138").Value = -16008.908

# ALC row 141
$ws_ALC.Range("H141").Value = 781.28
$ws_ALC.Range("I141").Value = 730.5
$ws_ALC.Range("K141").Value = 2191.5
$ws_ALC.Range("M141").Value = 2988.5

# ARM row 63
$ws_ARM.Range("H63").Value = 1798.2
$ws_ARM.Range("I63").Value = 1832.1111
$ws_ARM.Range("J63").Value = 1493
$ws_ARM.Range("K63").Value = 1832.1111
$ws_ARM.Range("L63").Value = 1493
$ws_ARM.Range("M63").Value = -1146.1111
$ws_ARM.Range("N63").Value = -2865

# ARM row 66
$ws_ARM.Range("H66").Value = 1798.2
$ws_ARM.Range("I66").Value = 1832.1111
$ws_ARM.Range("J66").Value = 1493
$ws_ARM.Range("K66").Value = 9160.5555
$ws_ARM.Range("L66").Value = 7465
$ws_ARM.Range("M66").Value = -5728.5555
$ws_ARM.Range("N66").Value = -14329

# ARM row 74
$ws_ARM.Range("H74").Value = 4632178.5
$ws_ARM.Range("I74").Value = 5953852.5
$ws_ARM.Range("K74").Value = 5953852.5
$ws_ARM.Range("M74").Value = -5952978.5

# ARM row 77
$ws_ARM.Range("H77").Value = 4632178.5
$ws_ARM.Range("I77").Value = 5953852.5
$ws_ARM.Range("K77").Value = 29769262.5
$ws_ARM.Range("M77").Value = -29764894.5

# ARM row 122
$ws_ARM.Range("H122").Value = 3152.875
$ws_ARM.Range("I122").Value = 3152.875
$ws_ARM.Range("K122").Value = 9458.625
$ws_ARM.Range("M122").Value = -7008.625

# ARM row 132
$ws_ARM.Range("H132").Value = 973068.8
$ws_ARM.Range("I132").Value = 1061119.2
$ws_ARM.Range("J132").Value = 4513.5
$ws_ARM.Range("K132").Value = 3183357.6
$ws_ARM.Range("L132").Value = 13540.5
$ws_ARM.Range("M132").Value = -3180827.6
$ws_ARM.Range("N132").Value = -18600.5

# ARM row 134
$ws_ARM.Range("H134").Value = 57526
$ws_ARM.Range("J134").Value = 57526
$ws_ARM.Range("L134").Value = 57526
$ws_ARM.Range("N134").Value = -67666

# BSM row 86
$ws_BSM.Range("H86").Value = 1508.75
$ws_BSM.Range("I86").Value = 1367.1428
$ws_BSM.Range("J86").Value = 2500
$ws_BSM.Range("K86").Value = 1367.1428
$ws_BSM.Range("L86").Value = 2500
$ws_BSM.Range("M86").Value = -244.1428000000001
$ws_BSM.Range("N86").Value = -4746

# BSM row 89
$ws_BSM.Range("H89").Value = 1508.75
$ws_BSM.Range("I89").Value = 1367.1428
$ws_BSM.Range("J89").Value = 2500
$ws_BSM.Range("K89").Value = 6835.714
$ws_BSM.Range("L89").Value = 12500
$ws_BSM.Range("M89").Value = -1219.714
$ws_BSM.Range("N89").Value = -23732

# BSM row 112
$ws_BSM.Range("H112").Value = 0
$ws_BSM.Range("J112").Value = 0
$ws_BSM.Range("L112").Value = 0
$ws_BSM.Range("N112").ClearContents() | Out-Null

# BSM row 134
$ws_BSM.Range("H134").Value = 420249.6
$ws_BSM.Range("I134").Value = 568612.8
$ws_BSM.Range("J134").Value = 4832.6665
$ws_BSM.Range("K134").Value = 1705838.4
$ws_BSM.Range("L134").Value = 14497.9995
$ws_BSM.Range("M134").Value = -1703303.4
$ws_BSM.Range("N134").Value = -19567.9995

# CRP row 31
$ws_CRP.Range("H31").Value = 115889.5
$ws_CRP.Range("I31").Value = 148568.5
$ws_CRP.Range("J31").Value = 43995.7
$ws_CRP.Range("K31").Value = 148568.5
$ws_CRP.Range("L31").Value = 43995.7
$ws_CRP.Range("M31").Value = -148273.5
$ws_CRP.Range("N31").Value = -44585.7

# CRP row 34
$ws_CRP.Range("H34").Value = 115889.5
$ws_CRP.Range("I34").Value = 148568.5
$ws_CRP.Range("J34").Value = 43995.7
$ws_CRP.Range("K34").Value = 148568.5
$ws_CRP.Range("L34").Value = 43995.7
$ws_CRP.Range("M34").Value = -148366.5
$ws_CRP.Range("N34").Value = -44399.7

# CRP row 58
$ws_CRP.Range("H58").Value = 303586.62
$ws_CRP.Range("I58").Value = 442562.7
$ws_CRP.Range("J58").Value = 4253.5386
$ws_CRP.Range("K58").Value = 442562.7
$ws_CRP.Range("L58").Value = 4253.5386
$ws_CRP.Range("M58").Value = -442359.7
$ws_CRP.Range("N58").Value = -4659.5386

# CRP row 99
$ws_CRP.Range("H99").Value = 3559
$ws_CRP.Range("I99").Value = 3036.2856
$ws_CRP.Range("J99").Value = 4473.75
$ws_CRP.Range("K99").Value = 3036.2856
$ws_CRP.Range("L99").Value = 4473.75
$ws_CRP.Range("M99").Value = -1538.2856
$ws_CRP.Range("N99").Value = -7469.75

# CRP row 126
$ws_CRP.Range("H126").Value = 3559
$ws_CRP.Range("I126").Value = 3036.2856
$ws_CRP.Range("J126").Value = 4473.75
$ws_CRP.Range("K126").Value = 9108.856800000001
$ws_CRP.Range("L126").Value = 13421.25
$ws_CRP.Range("M126").Value = -6638.856800000001
$ws_CRP.Range("N126").Value = -18361.25

# CRP row 134
$ws_CRP.Range("H134").Value = 13891.36
$ws_CRP.Range("I134").Value = 17890.264
$ws_CRP.Range("K134").Value = 53670.792
$ws_CRP.Range("M134").Value = -51135.792

# CRP row 136
$ws_CRP.Range("H136").Value = 303586.62
$ws_CRP.Range("I136").Value = 442562.7
$ws_CRP.Range("J136").Value = 4253.5386
$ws_CRP.Range("K136").Value = 1327688.1
$ws_CRP.Range("L136").Value = 12760.6158
$ws_CRP.Range("M136").Value = -1325138.1
$ws_CRP.Range("N136").Value = -17860.6158

# CUL row 122
$ws_CUL.Range("H122").Value = 21389332
$ws_CUL.Range("J122").Value = 42778144
$ws_CUL.Range("L122").Value = 385003296
$ws_CUL.Range("N122").Value = -385008196

# CUL row 131
$ws_CUL.Range("H131").Value = 26328.111
$ws_CUL.Range("I131").Value = 800
$ws_CUL.Range("J131").Value = 29519.125
$ws_CUL.Range("K131").Value = 2400
$ws_CUL.Range("L131").Value = 88557.375
$ws_CUL.Range("M131").Value = 2640
$ws_CUL.Range("N131").Value = -98637.375

# CUL row 140
$ws_CUL.Range("H140").Value = 3405.923
$ws_CUL.Range("I140").Value = 2660.0625
$ws_CUL.Range("K140").Value = 7980.1875
$ws_CUL.Range("M140").Value = -2800.1875

# CUL row 141
$ws_CUL.Range("H141").Value = 7038
$ws_CUL.Range("I141").Value = 7038
$ws_CUL.Range("K141").Value = 21114
$ws_CUL.Range("M141").Value = -15934

# GSM row 80
$ws_GSM.Range("H80").Value = 4607.067
$ws_GSM.Range("I80").Value = 3557.182
$ws_GSM.Range("K80").Value = 3557.182
$ws_GSM.Range("M80").Value = -2559.182

# GSM row 83
$ws_GSM.Range("H83").Value = 4607.067
$ws_GSM.Range("I83").Value = 3557.182
$ws_GSM.Range("K83").Value = 17785.91
$ws_GSM.Range("M83").Value = -12793.91

# GSM row 122
$ws_GSM.Range("H122").Value = 11862.25
$ws_GSM.Range("I122").Value = 5780
$ws_GSM.Range("J122").Value = 21999.334
$ws_GSM.Range("K122").Value = 17340
$ws_GSM.Range("L122").Value = 65998.00199999999
$ws_GSM.Range("M122").Value = -14890
$ws_GSM.Range("N122").Value = -70898.00199999999

# GSM row 132
$ws_GSM.Range("H132").Value = 256299.89
$ws_GSM.Range("I132").Value = 269897.34
$ws_GSM.Range("K132").Value = 809692.02
$ws_GSM.Range("M132").Value = -807162.02

# LTW row 45
$ws_LTW.Range("H45").Value = 0
$ws_LTW.Range("I45").Value = 0
$ws_LTW.Range("J45").Value = 0
$ws_LTW.Range("K45").Value = 0
$ws_LTW.Range("L45").Value = 0
$ws_LTW.Range("M45").ClearContents() | Out-Null
$ws_LTW.Range("N45").ClearContents() | Out-Null

# LTW row 136
$ws_LTW.Range("H136").Value = 37794.6
$ws_LTW.Range("I136").Value = 2186.2
$ws_LTW.Range("J136").Value = 251445
$ws_LTW.Range("K136").Value = 6558.599999999999
$ws_LTW.Range("L136").Value = 754335
$ws_LTW.Range("M136").Value = -4008.599999999999
$ws_LTW.Range("N136").Value = -759435

# WVR row 100
$ws_WVR.Range("H100").Value = 796.2917
$ws_WVR.Range("I100").Value = 755.95
$ws_WVR.Range("K100").Value = 1511.9
$ws_WVR.Range("M100").Value = -970.9000000000001

# WVR row 122
$ws_WVR.Range("H122").Value = 1709.7407
$ws_WVR.Range("I122").Value = 1380.2273
$ws_WVR.Range("J122").Value = 3159.6
$ws_WVR.Range("K122").Value = 4140.6819
$ws_WVR.Range("L122").Value = 9478.799999999999
$ws_WVR.Range("M122").Value = -1690.6819
$ws_WVR.Range("N122").Value = -14378.8

# WVR row 136
$ws_WVR.Range("H136").Value = 10496857
$ws_WVR.Range("I136").Value = 13211112
$ws_WVR.Range("K136").Value = 39633336
$ws_WVR.Range("M136").Value = -39630786
